# Regenerate orders with updated distance/size codes.
# The workbook stores trial "Condition", "Filename_Left", "Filename_Right",
# "Distance" and "Size" values that encode the viewing distance (Dxx) and
# image size (Sxx) used for each trial. This run updates those codes:
#   D64 -> D69
#   D80 -> D86
#   D51 -> D55
#   S30 -> S31
# across every cell on the sheet (headers are untouched since they don't
# contain these substrings).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Replace("D64", "D69")
$ws.Cells.Replace("D80", "D86")
$ws.Cells.Replace("D51", "D55")
$ws.Cells.Replace("S30", "S31")
